$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.302.66"
$ws.Range("E2").Value = "  -1.35%  "

# Row 3
$ws.Range("D3").Value = "2.521.67"
$ws.Range("E3").Value = "  -0.79%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.98"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = "  +3.10%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.78"
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = "  -5.13%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.572"
$ws.Range("D7").Style = $ws.Range("B7").Style
$ws.Range("E7").Value = "  -0.94%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("D9").Style = $ws.Range("B9").Style
$ws.Range("E9").Value = "  -3.13%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.40"
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value = "  -4.18%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0803"
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = "  -2.51%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.58"
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").Value = "  -0.64%  "

# Row 13
$ws.Range("E13").Value = "  -0.42%  "

# Row 14
$ws.Range("D14").Value = "2.907.32"
$ws.Range("E14").Value = "  -0.84%  "

# Row 15
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.44"
$ws.Range("D15").Style = $ws.Range("B15").Style
$ws.Range("E15").Value = "  +1.25%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.518.94"
$ws.Range("E16").Value = "  -2.10%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.849"
$ws.Range("D17").Style = $ws.Range("B17").Style
$ws.Range("E17").Value = "  -3.15%  "

# Row 18
$ws.Range("D18").Value = "42.401.38"
$ws.Range("E18").Value = "  -1.14%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.82"
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Value = "  -2.75%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.57"
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").Value = "  +0.24%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0956"
$ws.Range("E21").Value = "  -3.46%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.49"
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value = "  -1.70%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "249.36"
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value = "  -2.02%  "

# Row 24
$ws.Range("E24").Value = "  -0.58%  "

# Row 25
$ws.Range("E25").Value = "  -3.21%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.53"
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Value = "  -4.26%  "

# Row 27
$ws.Range("E27").Value = "  -0.21%  "

# Row 28
$ws.Range("E28").Value = "  +2.91%  "

# Row 29
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.09"
$ws.Range("D29").Style = $ws.Range("B29").Style
$ws.Range("E29").Value = "  -1.12%  "

# Row 30
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.72"
$ws.Range("D30").Style = $ws.Range("B30").Style
$ws.Range("E30").Value = "  +0.24%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.89"
$ws.Range("D31").Style = $ws.Range("B31").Style
$ws.Range("E31").Value = "  -4.46%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.78"
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").Value = "  -1.36%  "

# Row 33
$ws.Range("E33").Value = "  -1.23%  "

# Row 34
$ws.Range("E34").Value = "  +0.17%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.80"
$ws.Range("D35").Style = $ws.Range("B35").Style
$ws.Range("E35").Value = "  -0.97%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0778"
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").Value = "  -2.96%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.62"
$ws.Range("D37").Style = $ws.Range("B37").Style
$ws.Range("E37").Value = "  -0.87%  "

# Row 38
$ws.Range("E38").Value = "  -4.26%  "

# Row 39
$ws.Range("E39").Value = "  -1.66%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.74"
$ws.Range("D40").Style = $ws.Range("B40").Style
$ws.Range("E40").Value = "  -2.30%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.33"
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Value = "  +10.53%  "

# Row 42
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E42").Value = "  +0.38%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.79"
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Value = "  -2.80%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.30"
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = "  -5.02%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0298"
$ws.Range("D45").Style = $ws.Range("B45").Style
$ws.Range("E45").Value = "  -2.15%  "

# Row 46
$ws.Range("D46").Value = "2.014.25"
$ws.Range("E46").Value = "  -3.17%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.15"
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").Value = "  -2.65%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.79"
$ws.Range("D48").Style = $ws.Range("B48").Style
$ws.Range("E48").Value = "  -2.77%  "

# Row 49
$ws.Range("D49").Value = "2.763.25"
$ws.Range("E49").Value = "  -0.93%  "

# Row 50
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "101.61"
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value = "  -1.74%  "

# Row 51
$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.14"
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Value = "  -2.08%  "
